$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.646.24'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '1.696.59'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.48'
$ws.Range('E5').Value = '  -0.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3921'
$ws.Range('E7').Value = '  -0.53%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4037'
$ws.Range('E8').Value = '  +0.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.518'
$ws.Range('E9').Value = '  -0.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.002'
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.80'
$ws.Range('E11').Value = '  -2.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08838'
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.404'
$ws.Range('E13').Value = '  +2.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.54'
$ws.Range('E14').Value = '  +1.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.125'
$ws.Range('E15').Value = '  +6.67%  '
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('D17').Value = '1.704.39'
$ws.Range('E17').Value = '  +0.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '99.45'
$ws.Range('E18').Value = '  -0.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07012'
$ws.Range('E19').Value = '  -0.60%  '
$ws.Range('E20').Value = '  +0.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.066'
$ws.Range('E21').Value = '  +3.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.005'
$ws.Range('E22').Value = '  +0.44%  '
$ws.Range('E23').Value = '  +3.87%  '
$ws.Range('D24').Value = '24.647.99'
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.126'
$ws.Range('E25').Value = '  +3.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.344'
$ws.Range('E26').Value = '  +1.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.62'
$ws.Range('E27').Value = '  +0.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '163.04'
$ws.Range('E28').Value = '  +1.98%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.717'
$ws.Range('E29').Value = '  +16.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '135.40'
$ws.Range('E30').Value = '  +0.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.147'
$ws.Range('E31').Value = '  -1.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08945'
$ws.Range('E32').Value = '  +4.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.588'
$ws.Range('E33').Value = '  +3.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.066'
$ws.Range('E34').Value = '  -4.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.964'
$ws.Range('E35').Value = '  +0.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '11.04'
$ws.Range('E36').Value = '  -3.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2743'
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02881'
$ws.Range('E38').Value = '  +4.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.41'
$ws.Range('E39').Value = '  -1.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09148'
$ws.Range('E40').Value = '  +0.80%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.456'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7643'
$ws.Range('E42').Value = '  -0.66%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '15.79'
$ws.Range('E43').Value = '  +2.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7161'
$ws.Range('E44').Value = '  -0.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.557'
$ws.Range('E45').Value = '  +1.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.214'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.001'
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.335'
$ws.Range('E48').Value = '  -1.46%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '139.87'
$ws.Range('E49').Value = '  -1.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07964'
$ws.Range('E50').Value = '  -0.74%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '90.27'
$ws.Range('E51').Value = '  +2.04%  '
